$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.171.31'
$ws.Range("E2").Value = '  -3.24%  '

$ws.Range("D3").Value = '2.389.54'
$ws.Range("E3").Value = '  +5.80%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '292.68'
$ws.Range("E5").Value = '  -2.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.58'
$ws.Range("E6").Value = '  -6.77%  '

$ws.Range("E7").Value = '  -0.91%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  -2.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.17'
$ws.Range("E10").Value = '  -3.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0775'
$ws.Range("E11").Value = '  -0.83%  '

$ws.Range("E12").Value = '  -2.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  +0.96%  '

$ws.Range("D14").Value = '2.754.63'
$ws.Range("E14").Value = '  +5.51%  '

$ws.Range("D15").Value = '2.386.30'
$ws.Range("E15").Value = '  +5.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.11'
$ws.Range("E16").Value = '  +4.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.824'
$ws.Range("E17").Value = '  +3.69%  '

$ws.Range("D18").Value = '45.166.30'
$ws.Range("E18").Value = '  -3.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.39'
$ws.Range("E19").Value = '  -4.12%  '

$ws.Range("D20").Value = '0.0₃0932'
$ws.Range("E20").Value = '  +0.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.08'
$ws.Range("E21").Value = '  +3.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.41'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.58'
$ws.Range("E23").Value = '  -4.77%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.75'
$ws.Range("E24").Value = '  -2.96%  '

$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("E26").Value = '  +0.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.22'
$ws.Range("E27").Value = '  -0.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.26'
$ws.Range("E28").Value = '  -13.21%  '

$ws.Range("E29").Value = '  -2.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.80'
$ws.Range("E30").Value = '  +19.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.88'
$ws.Range("E31").Value = '  +5.26%  '

$ws.Range("E32").Value = '  -3.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '146.77'
$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.38'
$ws.Range("E34").Value = '  -0.92%  '

$ws.Range("E35").Value = '  -1.77%  '

$ws.Range("E36").Value = '  +14.09%  '

$ws.Range("E37").Value = '  -1.86%  '

$ws.Range("E38").Value = '  -1.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.43'
$ws.Range("E39").Value = '  -11.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.67'
$ws.Range("E40").Value = '  -4.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0291'
$ws.Range("E41").Value = '  -1.79%  '

$ws.Range("D42").Value = '1.970.18'
$ws.Range("E42").Value = '  +8.81%  '

$ws.Range("E43").Value = '  -1.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '87.99'
$ws.Range("E45").Value = '  -2.90%  '

$ws.Range("E46").Value = '  -14.21%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.28'
$ws.Range("E47").Value = '  +19.81%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.35'
$ws.Range("E48").Value = '  +7.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '99.04'
$ws.Range("E49").Value = '  +5.72%  '

$ws.Range("D50").Value = '2.625.27'
$ws.Range("E50").Value = '  +5.53%  '

$ws.Range("E51").Value = '  -3.48%  '
